$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 110; $r++) {
    $fVal = $ws.Cells.Item($r, 6).Value2
    $gVal = $ws.Cells.Item($r, 7).Value2
    $hVal = $ws.Cells.Item($r, 8).Value2
    $iVal = $ws.Cells.Item($r, 9).Value2

    $ws.Cells.Item($r, 6).Value2 = $gVal
    $ws.Cells.Item($r, 7).Value2 = $fVal
    $ws.Cells.Item($r, 8).Value2 = $iVal
    $ws.Cells.Item($r, 9).Value2 = $hVal
}
